$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.522.09"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.280.67"
$ws.Range("E3").Value = "  +3.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "616.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.283.29"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.498"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "3.804.72"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "66.498.72"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "3.268.43"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "507.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.759"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.131"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +49.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +21.08%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0796"
$ws.Range("E38").Value = "  +16.56%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "495.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0426"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.59%  "
$ws.Range("D45").Value = "3.028.91"
$ws.Range("E45").Value = "  +7.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.295"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.27%  "
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.77%  "
